# Update attendance/visitor numbers in the "展览" (rId1/sheet1) and
# "全部类型" (rId4/sheet4) worksheets, as produced by the site's data
# regeneration (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F4").Value = 13566
    $ws.Range("F17").Value = 439
    $ws.Range("F21").Value = 967
}
